# Updated cryptos list on Sat Mar 16 20:57:59 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# re-sorts two pairs of rows whose rank flipped (Dogecoin/Avalanche and
# Chainlink/WrappedBTC) by swapping their Coin/Link/Price/Volume cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    # Force the literal string into the cell even when it looks numeric
    # (e.g. "1.00" or "587.92"), then restore the cell to the workbook's
    # default ("Normal") style so no stray number-format style is left
    # behind, matching how the sheet was originally authored.
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

# Row -> (newPrice, newVolume) for rows whose Coin/Link didn't move.
$priceUpdates = @(
    @{ Row = 2;  Price = "67.192.19";  Volume = "  -1.11%  " },
    @{ Row = 3;  Price = "3.619.44";   Volume = "  -0.25%  " },
    @{ Row = 4;  Price = "1.00";       Volume = "  +0.48%  " },
    @{ Row = 5;  Price = "587.92";     Volume = "  +0.48%  " },
    @{ Row = 6;  Price = "184.02";     Volume = "  +4.49%  " },
    @{ Row = 7;  Price = "0.612";      Volume = "  -1.82%  " },
    @{ Row = 8;  Price = "1.00";       Volume = "  -0.01%  " },
    @{ Row = 9;  Price = "0.671";      Volume = "  -4.41%  " },
    @{ Row = 12; Price = "0.0000253";  Volume = "  -11.51%  " },
    @{ Row = 13; Price = "9.91";       Volume = "  -5.49%  " },
    @{ Row = 14; Price = "4.198.53";   Volume = "  +0.08%  " },
    @{ Row = 15; Price = "3.620.32";   Volume = "  -0.20%  " },
    @{ Row = 19; Price = "12.20";      Volume = "  -3.16%  " },
    @{ Row = 21; Price = "391.85";     Volume = "  -3.23%  " },
    @{ Row = 22; Price = "4.31";       Volume = "  -3.85%  " },
    @{ Row = 23; Price = "84.91";      Volume = "  -2.99%  " },
    @{ Row = 24; Price = "2.87";       Volume = "  -3.67%  " },
    @{ Row = 25; Price = "12.25";      Volume = "  -2.34%  " },
    @{ Row = 26; Price = "6.06";       Volume = "  +0.79%  " },
    @{ Row = 27; Price = "10.29";      Volume = "  -2.99%  " },
    @{ Row = 28; Price = "3.60";       Volume = "  -10.78%  " },
    @{ Row = 29; Price = "8.97";       Volume = "  -3.96%  " },
    @{ Row = 30; Price = "31.11";      Volume = "  -3.49%  " },
    @{ Row = 31; Price = "6.79";       Volume = "  -4.71%  " },
    @{ Row = 32; Price = "65.46";      Volume = "  +2.11%  " },
    @{ Row = 33; Price = "11.88";      Volume = "  -2.62%  " },
    @{ Row = 34; Price = "594.55";     Volume = "  +0.57%  " },
    @{ Row = 35; Price = "0.112";      Volume = "  -2.70%  " },
    @{ Row = 36; Price = "41.50";      Volume = "  -2.29%  " },
    @{ Row = 38; Price = "1.00";       Volume = "  +0.55%  " },
    @{ Row = 39; Price = "0.373";      Volume = "  -4.97%  " },
    @{ Row = 40; Price = "0.0₃0739";   Volume = "  -15.42%  " },
    @{ Row = 42; Price = "2.77";       Volume = "  -6.92%  " },
    @{ Row = 43; Price = "0.0411";     Volume = "  -4.66%  " },
    @{ Row = 44; Price = "2.40";       Volume = "  -9.87%  " },
    @{ Row = 45; Price = "2.701.85";   Volume = "  +0.96%  " },
    @{ Row = 46; Price = "0.129";      Volume = "  -2.60%  " },
    @{ Row = 47; Price = "3.04";       Volume = "  -1.56%  " },
    @{ Row = 49; Price = "135.62";     Volume = "  -3.02%  " },
    @{ Row = 50; Price = "8.25";       Volume = "  -7.04%  " },
    @{ Row = 51; Price = "2.57";       Volume = "  -5.41%  " }
)

foreach ($u in $priceUpdates) {
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.Price
    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}

# Rows where only the Volume(1h) figure moved (Price unchanged).
$volumeOnlyUpdates = @(
    @{ Row = 16; Volume = "  -0.03%  " },
    @{ Row = 20; Volume = "  -4.06%  " },
    @{ Row = 37; Volume = "  +0.33%  " },
    @{ Row = 41; Volume = "  -1.51%  " },
    @{ Row = 48; Volume = "  -5.44%  " }
)

foreach ($u in $volumeOnlyUpdates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.Volume
}

# Rows 10/11 and 17/18 swapped rank - update Coin / Link / Price / Volume
# for each row so the pair trades places.
$ws.Cells.Item(10, 2).Value = "Avalanche"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Cells.Item(10, 4) "53.78"
$ws.Cells.Item(10, 5).Value = "  -2.29%  "

$ws.Cells.Item(11, 2).Value = "Dogecoin"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Cells.Item(11, 4) "0.144"
$ws.Cells.Item(11, 5).Value = "  -9.16%  "

$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Cells.Item(17, 4) "67.258.37"
$ws.Cells.Item(17, 5).Value = "  -0.36%  "

$ws.Cells.Item(18, 2).Value = "Chainlink"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(18, 4) "18.35"
$ws.Cells.Item(18, 5).Value = "  -3.84%  "
